$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold font, border, centered/top alignment) from the row
# above so the new row's label cell (A16) matches the existing style used
# by A2:A15.
$ws.Range("A15").Copy($ws.Range("A16"))

# Row 16 data: HKL index 14, "HexGrid-60degTilt5degRes" dataset
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 0.9842483788308727
$ws.Range("D16").Value = 1.021358771308097
$ws.Range("E16").Value = 0.9860784505992501
$ws.Range("F16").Value = 0.9842483788308727
$ws.Range("G16").Value = 1.01029065399248
$ws.Range("H16").Value = 0.9735217150096211
$ws.Range("I16").Value = 0.9831372741286619
$ws.Range("J16").Value = 1.021358771308097
$ws.Range("K16").Value = 1.003718610953674
$ws.Range("L16").Value = 0.993983494892273
$ws.Range("M16").Value = 0.9931058739781639
